$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)

# Row 9
$ws.Range("H9").Value = 5145.048
$ws.Range("I9").Value = 6691.4375
$ws.Range("J9").Value = 196.6
$ws.Range("K9").Value = 6691.4375
$ws.Range("L9").Value = 196.6
$ws.Range("M9").Value = -6522.4375
$ws.Range("N9").Value = -534.6

# Row 12
$ws.Range("H12").Value = 8507.583
$ws.Range("I12").Value = 12650.125
$ws.Range("K12").Value = 12650.125
$ws.Range("M12").Value = -12480.125

# Row 18
$ws.Range("H18").Value = 2000
$ws.Range("I18").Value = 2000
$ws.Range("K18").Value = 2000
$ws.Range("M18").Value = -1716

# Row 21
$ws.Range("H21").Value = 5999
$ws.Range("I21").Value = 5999
$ws.Range("K21").Value = 5999
$ws.Range("M21").Value = -5531

# Row 23
$ws.Range("H23").Value = 5999
$ws.Range("I23").Value = 5999
$ws.Range("K23").Value = 5999
$ws.Range("M23").Value = -5765

# Row 38
$ws.Range("H38").Value = 1766.1111
$ws.Range("I38").Value = 1486.875
$ws.Range("K38").Value = 4460.625
$ws.Range("M38").Value = -4088.625

# Row 43
$ws.Range("H43").Value = 4599.8
$ws.Range("I43").Value = 3999.75
$ws.Range("J43").Value = 7000
$ws.Range("K43").Value = 3999.75
$ws.Range("L43").Value = 7000
$ws.Range("M43").Value = -3930.75
$ws.Range("N43").Value = -7138

# Row 70
$ws.Range("H70").Value = 7091.864
$ws.Range("J70").Value = 8961
$ws.Range("L70").Value = 26883
$ws.Range("N70").Value = -27423

# Row 73
$ws.Range("H73").Value = 7091.864
$ws.Range("J73").Value = 8961
$ws.Range("L73").Value = 26883
$ws.Range("N73").Value = -28755

# Row 113
$ws.Range("H113").Value = 1557
$ws.Range("I113").Value = 1564.125
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1564.125
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1689.875
$ws.Range("N113").Value = -8008

# Row 132
$ws.Range("H132").Value = 1465.75
$ws.Range("I132").Value = 1461.4688
$ws.Range("K132").Value = 4384.4064
$ws.Range("M132").Value = -1854.4064

# Row 137
$ws.Range("H137").Value = 4335.448
$ws.Range("I137").Value = 4099.522
$ws.Range("J137").Value = 5239.8335
$ws.Range("K137").Value = 12298.566
$ws.Range("L137").Value = 15719.5005
$ws.Range("M137").Value = -9748.565999999999
$ws.Range("N137").Value = -20819.5005

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)

# Row 6
$ws.Range("H6").Value = 17833
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 17833
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 17833
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -18179

# Row 10
$ws.Range("H10").Value = 11665.667
$ws.Range("I10").Value = 15000
$ws.Range("J10").Value = 9998.5
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 9998.5
$ws.Range("M10").Value = -14830
$ws.Range("N10").Value = -10338.5

# Row 11
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# Row 45
$ws.Range("H45").Value = 1714.7693
$ws.Range("I45").Value = 1493.8334
$ws.Range("J45").Value = 1904.1428
$ws.Range("K45").Value = 1493.8334
$ws.Range("L45").Value = 1904.1428
$ws.Range("M45").Value = -1116.8334
$ws.Range("N45").Value = -2658.1428

# Row 74
$ws.Range("H74").Value = 125007660
$ws.Range("I74").Value = 200005260
$ws.Range("J74").Value = 11666.667
$ws.Range("K74").Value = 200005260
$ws.Range("L74").Value = 11666.667
$ws.Range("M74").Value = -200004386
$ws.Range("N74").Value = -13414.667

# Row 77
$ws.Range("H77").Value = 125007660
$ws.Range("I77").Value = 200005260
$ws.Range("J77").Value = 11666.667
$ws.Range("K77").Value = 1000026300
$ws.Range("L77").Value = 58333.335
$ws.Range("M77").Value = -1000021932
$ws.Range("N77").Value = -67069.33499999999

# Row 132
$ws.Range("H132").Value = 4753824
$ws.Range("I132").Value = 2568024.5
$ws.Range("K132").Value = 7704073.5
$ws.Range("M132").Value = -7701543.5

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)

# Row 4
$ws.Range("H4").Value = 492.98386
$ws.Range("I4").Value = 487.73334
$ws.Range("J4").Value = 506.88235
$ws.Range("K4").Value = 487.73334
$ws.Range("L4").Value = 506.88235
$ws.Range("M4").Value = -372.73334
$ws.Range("N4").Value = -736.88235

# Row 86
$ws.Range("H86").Value = 2961.12
$ws.Range("I86").Value = 2891.3157
$ws.Range("J86").Value = 3182.1667
$ws.Range("K86").Value = 2891.3157
$ws.Range("L86").Value = 3182.1667
$ws.Range("M86").Value = -1768.3157
$ws.Range("N86").Value = -5428.1667

# Row 89
$ws.Range("H89").Value = 2961.12
$ws.Range("I89").Value = 2891.3157
$ws.Range("J89").Value = 3182.1667
$ws.Range("K89").Value = 14456.5785
$ws.Range("L89").Value = 15910.8335
$ws.Range("M89").Value = -8840.5785
$ws.Range("N89").Value = -27142.8335

# Row 134
$ws.Range("H134").Value = 17785250
$ws.Range("I134").Value = 21468970
$ws.Range("K134").Value = 64406910
$ws.Range("M134").Value = -64404375

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)

# Row 21
$ws.Range("H21").Value = 14450
$ws.Range("I21").Value = 12000
$ws.Range("J21").Value = 16900
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 16900
$ws.Range("M21").Value = -11765
$ws.Range("N21").Value = -17370

# Row 31
$ws.Range("H31").Value = 9256.789
$ws.Range("I31").Value = 8557.8
$ws.Range("J31").Value = 9506.429
$ws.Range("K31").Value = 8557.8
$ws.Range("L31").Value = 9506.429
$ws.Range("M31").Value = -8262.8
$ws.Range("N31").Value = -10096.429

# Row 34
$ws.Range("H34").Value = 9256.789
$ws.Range("I34").Value = 8557.8
$ws.Range("J34").Value = 9506.429
$ws.Range("K34").Value = 8557.8
$ws.Range("L34").Value = 9506.429
$ws.Range("M34").Value = -8355.8
$ws.Range("N34").Value = -9910.429

# Row 134
$ws.Range("H134").Value = 8930342
$ws.Range("I134").Value = 10001767
$ws.Range("K134").Value = 30005301
$ws.Range("M134").Value = -30002766

# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)

# Row 41
$ws.Range("H41").Value = 3248.875
$ws.Range("I41").Value = 2297.6
$ws.Range("J41").Value = 4834.3335
$ws.Range("K41").Value = 6892.799999999999
$ws.Range("L41").Value = 14503.0005
$ws.Range("M41").Value = -6554.799999999999
$ws.Range("N41").Value = -15179.0005

# Row 131
$ws.Range("H131").Value = 2700
$ws.Range("J131").Value = 3000
$ws.Range("L131").Value = 9000
$ws.Range("N131").Value = -19080

# Row 139
$ws.Range("H139").Value = 1588.4445
$ws.Range("I139").Value = 1588.4445
$ws.Range("K139").Value = 4765.333500000001
$ws.Range("M139").Value = 374.6664999999994

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)

# Row 132
$ws.Range("H132").Value = 2566605.2
$ws.Range("I132").Value = 2619034.5
$ws.Range("K132").Value = 7857103.5
$ws.Range("M132").Value = -7854573.5

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)

# Row 22
$ws.Range("H22").Value = 2804.3333
$ws.Range("I22").Value = 2430.3333
$ws.Range("K22").Value = 2430.3333
$ws.Range("M22").Value = -2135.3333

# Row 27
$ws.Range("H27").Value = 2804.3333
$ws.Range("I27").Value = 2430.3333
$ws.Range("K27").Value = 2430.3333
$ws.Range("M27").Value = -2323.3333

# Row 34
$ws.Range("H34").Value = 21
$ws.Range("I34").Value = 21
$ws.Range("K34").Value = 21
$ws.Range("M34").Value = 151

# Row 40
$ws.Range("H40").Value = 4224.25
$ws.Range("I40").Value = 4224.25
$ws.Range("K40").Value = 4224.25
$ws.Range("M40").Value = -4088.25

# Row 122
$ws.Range("H122").Value = 4798.4
$ws.Range("I122").Value = 4798.4
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14395.2
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11945.2
$ws.Range("N122").ClearContents()

# Row 132
$ws.Range("H132").Value = 9266092
$ws.Range("I132").Value = 11370831
$ws.Range("K132").Value = 34112493
$ws.Range("M132").Value = -34109963

# --- Sheet: WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)

# Row 122
$ws.Range("H122").Value = 3235.75
$ws.Range("I122").Value = 3355.4285
$ws.Range("K122").Value = 10066.2855
$ws.Range("M122").Value = -7616.2855
